$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.861094666666666
$ws.Range("H2").Value = 29.583284
$ws.Range("I2").Value = 0.243709096397741
$ws.Range("J2").Value = 0.2437090963977409
$ws.Range("M2").Value = 0.029424
$ws.Range("N2").Value = 0.08827199999999999
$ws.Range("O2").Value = 0.1473063425232919
$ws.Range("P2").Value = 0.1473063425232919
$ws.Range("Q2").Value = 0.2901528494719999
$ws.Range("R2").Value = 2.611375645248
$ws.Range("S2").Value = 0.0358998956300076
$ws.Range("T2").Value = 0.0358998956300076

# Row 3
$ws.Range("G3").Value = 9.861094666666666
$ws.Range("H3").Value = 29.583284
$ws.Range("I3").Value = 0.243709096397741
$ws.Range("J3").Value = 0.2437090963977409
$ws.Range("O3").Value = 0.852693657476708
$ws.Range("P3").Value = 0.852693657476708
$ws.Range("Q3").Value = 1.679571226910667
$ws.Range("R3").Value = 15.116141042196
$ws.Range("S3").Value = 0.2078092007677333
$ws.Range("T3").Value = 0.2078092007677333

# Row 4
$ws.Range("I4").Value = 0.7254466225154019
$ws.Range("J4").Value = 0.7254466225154018
$ws.Range("M4").Value = 0.029424
$ws.Range("N4").Value = 0.08827199999999999
$ws.Range("O4").Value = 0.1473063425232919
$ws.Range("P4").Value = 0.1473063425232919
$ws.Range("Q4").Value = 0.8636953145119999
$ws.Range("R4").Value = 7.773257830607998
$ws.Range("S4").Value = 0.1068628886586191
$ws.Range("T4").Value = 0.1068628886586191

# Row 5
$ws.Range("I5").Value = 0.7254466225154019
$ws.Range("J5").Value = 0.7254466225154018
$ws.Range("O5").Value = 0.852693657476708
$ws.Range("P5").Value = 0.852693657476708
$ws.Range("S5").Value = 0.6185837338567828
$ws.Range("T5").Value = 0.6185837338567827

# Row 6
$ws.Range("I6").Value = 0.03084428108685718
$ws.Range("J6").Value = 0.03084428108685716
$ws.Range("M6").Value = 0.029424
$ws.Range("N6").Value = 0.08827199999999999
$ws.Range("O6").Value = 0.1473063425232919
$ws.Range("P6").Value = 0.1473063425232919
$ws.Range("Q6").Value = 0.036722289728
$ws.Range("R6").Value = 0.3305006075519999
$ws.Range("S6").Value = 0.004543558234665279
$ws.Range("T6").Value = 0.004543558234665277

# Row 7
$ws.Range("I7").Value = 0.03084428108685718
$ws.Range("J7").Value = 0.03084428108685716
$ws.Range("O7").Value = 0.852693657476708
$ws.Range("P7").Value = 0.852693657476708
$ws.Range("S7").Value = 0.0263007228521919
$ws.Range("T7").Value = 0.02630072285219189
